$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-07-13 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-14 Friday", 2) | Out-Null
$d.Content.Find.Execute("44×60=", $true, $false, $false, $false, $false, $true, 1, $false, "52×29=", 2) | Out-Null
$d.Content.Find.Execute("70×45=", $true, $false, $false, $false, $false, $true, 1, $false, "18×75=", 2) | Out-Null
$d.Content.Find.Execute("32×85=", $true, $false, $false, $false, $false, $true, 1, $false, "89×24=", 2) | Out-Null
$d.Content.Find.Execute("87×33=", $true, $false, $false, $false, $false, $true, 1, $false, "100×96=", 2) | Out-Null
$d.Content.Find.Execute("64×38=", $true, $false, $false, $false, $false, $true, 1, $false, "60×92=", 2) | Out-Null
$d.Content.Find.Execute("60×95=", $true, $false, $false, $false, $false, $true, 1, $false, "14×61=", 2) | Out-Null
$d.Content.Find.Execute("15×79=", $true, $false, $false, $false, $false, $true, 1, $false, "77×65=", 2) | Out-Null
$d.Content.Find.Execute("17×21=", $true, $false, $false, $false, $false, $true, 1, $false, "38×17=", 2) | Out-Null
$d.Content.Find.Execute("34×21=", $true, $false, $false, $false, $false, $true, 1, $false, "15×40=", 2) | Out-Null
$d.Content.Find.Execute("35×47=", $true, $false, $false, $false, $false, $true, 1, $false, "73×86=", 2) | Out-Null
$d.Content.Find.Execute("13×40=", $true, $false, $false, $false, $false, $true, 1, $false, "81×68=", 2) | Out-Null
$d.Content.Find.Execute("20×32=", $true, $false, $false, $false, $false, $true, 1, $false, "71×72=", 2) | Out-Null
$d.Content.Find.Execute("53×80=", $true, $false, $false, $false, $false, $true, 1, $false, "94×66=", 2) | Out-Null
$d.Content.Find.Execute("37×94=", $true, $false, $false, $false, $false, $true, 1, $false, "37×63=", 2) | Out-Null
$d.Content.Find.Execute("87×73=", $true, $false, $false, $false, $false, $true, 1, $false, "50×70=", 2) | Out-Null
$d.Content.Find.Execute("22×76=", $true, $false, $false, $false, $false, $true, 1, $false, "65×60=", 2) | Out-Null
$d.Content.Find.Execute("56×72=", $true, $false, $false, $false, $false, $true, 1, $false, "38×98=", 2) | Out-Null
$d.Content.Find.Execute("76×74=", $true, $false, $false, $false, $false, $true, 1, $false, "81×21=", 2) | Out-Null
$d.Content.Find.Execute("69×28=", $true, $false, $false, $false, $false, $true, 1, $false, "56×59=", 2) | Out-Null
$d.Content.Find.Execute("70×59=", $true, $false, $false, $false, $false, $true, 1, $false, "95×85=", 2) | Out-Null
$d.Content.Find.Execute("43×67=", $true, $false, $false, $false, $false, $true, 1, $false, "31×24=", 2) | Out-Null
$d.Content.Find.Execute("73×53=", $true, $false, $false, $false, $false, $true, 1, $false, "44×55=", 2) | Out-Null
$d.Content.Find.Execute("13×44=", $true, $false, $false, $false, $false, $true, 1, $false, "57×48=", 2) | Out-Null
$d.Content.Find.Execute("23×44=", $true, $false, $false, $false, $false, $true, 1, $false, "33×89=", 2) | Out-Null
$d.Content.Find.Execute("10×13=", $true, $false, $false, $false, $false, $true, 1, $false, "98×68=", 2) | Out-Null
$d.Content.Find.Execute("98×57=", $true, $false, $false, $false, $false, $true, 1, $false, "18×63=", 2) | Out-Null
$d.Content.Find.Execute("83×56=", $true, $false, $false, $false, $false, $true, 1, $false, "11×10=", 2) | Out-Null
$d.Content.Find.Execute("17×32=", $true, $false, $false, $false, $false, $true, 1, $false, "52×55=", 2) | Out-Null
$d.Content.Find.Execute("79×76=", $true, $false, $false, $false, $false, $true, 1, $false, "96×69=", 2) | Out-Null
$d.Content.Find.Execute("11×86=", $true, $false, $false, $false, $false, $true, 1, $false, "68×38=", 2) | Out-Null
$d.Content.Find.Execute("43×80=", $true, $false, $false, $false, $false, $true, 1, $false, "21×38=", 2) | Out-Null
$d.Content.Find.Execute("64×28=", $true, $false, $false, $false, $false, $true, 1, $false, "43×81=", 2) | Out-Null
$d.Content.Find.Execute("32×58=", $true, $false, $false, $false, $false, $true, 1, $false, "67×61=", 2) | Out-Null
$d.Content.Find.Execute("92×74=", $true, $false, $false, $false, $false, $true, 1, $false, "26×27=", 2) | Out-Null
$d.Content.Find.Execute("12×71=", $true, $false, $false, $false, $false, $true, 1, $false, "60×61=", 2) | Out-Null
$d.Content.Find.Execute("23×81=", $true, $false, $false, $false, $false, $true, 1, $false, "27×53=", 2) | Out-Null
$d.Content.Find.Execute("50×88=", $true, $false, $false, $false, $false, $true, 1, $false, "13×90=", 2) | Out-Null
$d.Content.Find.Execute("74×69=", $true, $false, $false, $false, $false, $true, 1, $false, "21×76=", 2) | Out-Null
$d.Content.Find.Execute("76×13=", $true, $false, $false, $false, $false, $true, 1, $false, "60×26=", 2) | Out-Null
$d.Content.Find.Execute("27×20=", $true, $false, $false, $false, $false, $true, 1, $false, "92×33=", 2) | Out-Null
$d.Content.Find.Execute("31×12=", $true, $false, $false, $false, $false, $true, 1, $false, "44×52=", 2) | Out-Null
$d.Content.Find.Execute("89×65=", $true, $false, $false, $false, $false, $true, 1, $false, "97×57=", 2) | Out-Null
$d.Content.Find.Execute("33×40=", $true, $false, $false, $false, $false, $true, 1, $false, "97×19=", 2) | Out-Null
$d.Content.Find.Execute("50×21=", $true, $false, $false, $false, $false, $true, 1, $false, "99×23=", 2) | Out-Null
$d.Content.Find.Execute("71×54=", $true, $false, $false, $false, $false, $true, 1, $false, "54×13=", 2) | Out-Null
$d.Content.Find.Execute("78×16=", $true, $false, $false, $false, $false, $true, 1, $false, "82×75=", 2) | Out-Null
$d.Content.Find.Execute("44×70=", $true, $false, $false, $false, $false, $true, 1, $false, "24×94=", 2) | Out-Null
$d.Content.Find.Execute("83×76=", $true, $false, $false, $false, $false, $true, 1, $false, "11×58=", 2) | Out-Null
$d.Content.Find.Execute("63×74=", $true, $false, $false, $false, $false, $true, 1, $false, "25×88=", 2) | Out-Null
$d.Content.Find.Execute("97×23=", $true, $false, $false, $false, $false, $true, 1, $false, "10×96=", 2) | Out-Null
$d.Content.Find.Execute("47×62=", $true, $false, $false, $false, $false, $true, 1, $false, "91×89=", 2) | Out-Null
$d.Content.Find.Execute("99×80=", $true, $false, $false, $false, $false, $true, 1, $false, "62×94=", 2) | Out-Null
$d.Content.Find.Execute("57×90=", $true, $false, $false, $false, $false, $true, 1, $false, "30×42=", 2) | Out-Null
$d.Content.Find.Execute("10×46=", $true, $false, $false, $false, $false, $true, 1, $false, "11×56=", 2) | Out-Null
$d.Content.Find.Execute("90×23=", $true, $false, $false, $false, $false, $true, 1, $false, "71×65=", 2) | Out-Null
$d.Content.Find.Execute("16×27=", $true, $false, $false, $false, $false, $true, 1, $false, "17×28=", 2) | Out-Null
$d.Content.Find.Execute("52×39=", $true, $false, $false, $false, $false, $true, 1, $false, "59×66=", 2) | Out-Null
$d.Content.Find.Execute("87×52=", $true, $false, $false, $false, $false, $true, 1, $false, "29×22=", 2) | Out-Null
$d.Content.Find.Execute("15×71=", $true, $false, $false, $false, $false, $true, 1, $false, "56×47=", 2) | Out-Null
$d.Content.Find.Execute("37×12=", $true, $false, $false, $false, $false, $true, 1, $false, "96×75=", 2) | Out-Null
$d.Content.Find.Execute("90×40=", $true, $false, $false, $false, $false, $true, 1, $false, "57×43=", 2) | Out-Null
$d.Content.Find.Execute("42×43=", $true, $false, $false, $false, $false, $true, 1, $false, "16×73=", 2) | Out-Null
$d.Content.Find.Execute("30×60=", $true, $false, $false, $false, $false, $true, 1, $false, "19×57=", 2) | Out-Null
$d.Content.Find.Execute("84×75=", $true, $false, $false, $false, $false, $true, 1, $false, "88×35=", 2) | Out-Null
$d.Content.Find.Execute("22×61=", $true, $false, $false, $false, $false, $true, 1, $false, "62×44=", 2) | Out-Null
$d.Content.Find.Execute("16×17=", $true, $false, $false, $false, $false, $true, 1, $false, "21×75=", 2) | Out-Null
$d.Content.Find.Execute("83×77=", $true, $false, $false, $false, $false, $true, 1, $false, "51×65=", 2) | Out-Null
$d.Content.Find.Execute("78×64=", $true, $false, $false, $false, $false, $true, 1, $false, "13×78=", 2) | Out-Null
$d.Content.Find.Execute("24×48=", $true, $false, $false, $false, $false, $true, 1, $false, "22×65=", 2) | Out-Null
$d.Content.Find.Execute("100×35=", $true, $false, $false, $false, $false, $true, 1, $false, "46×39=", 2) | Out-Null
$d.Content.Find.Execute("14×40=", $true, $false, $false, $false, $false, $true, 1, $false, "15×61=", 2) | Out-Null
$d.Content.Find.Execute("84×98=", $true, $false, $false, $false, $false, $true, 1, $false, "20×81=", 2) | Out-Null
$d.Content.Find.Execute("45×29=", $true, $false, $false, $false, $false, $true, 1, $false, "65×34=", 2) | Out-Null
$d.Content.Find.Execute("85×15=", $true, $false, $false, $false, $false, $true, 1, $false, "13×19=", 2) | Out-Null
$d.Content.Find.Execute("72×80=", $true, $false, $false, $false, $false, $true, 1, $false, "45×36=", 2) | Out-Null
$d.Content.Find.Execute("53×29=", $true, $false, $false, $false, $false, $true, 1, $false, "57×73=", 2) | Out-Null
$d.Content.Find.Execute("92×42=", $true, $false, $false, $false, $false, $true, 1, $false, "35×75=", 2) | Out-Null
$d.Content.Find.Execute("19×75=", $true, $false, $false, $false, $false, $true, 1, $false, "93×29=", 2) | Out-Null
$d.Content.Find.Execute("36×49=", $true, $false, $false, $false, $false, $true, 1, $false, "11×51=", 2) | Out-Null
$d.Content.Find.Execute("54×33=", $true, $false, $false, $false, $false, $true, 1, $false, "22×49=", 2) | Out-Null
$d.Content.Find.Execute("12×73=", $true, $false, $false, $false, $false, $true, 1, $false, "28×78=", 2) | Out-Null
$d.Content.Find.Execute("14×60=", $true, $false, $false, $false, $false, $true, 1, $false, "18×29=", 2) | Out-Null
$d.Content.Find.Execute("99×33=", $true, $false, $false, $false, $false, $true, 1, $false, "78×86=", 2) | Out-Null
$d.Content.Find.Execute("58×76=", $true, $false, $false, $false, $false, $true, 1, $false, "20×60=", 2) | Out-Null
$d.Content.Find.Execute("92×73=", $true, $false, $false, $false, $false, $true, 1, $false, "45×68=", 2) | Out-Null
$d.Content.Find.Execute("98×45=", $true, $false, $false, $false, $false, $true, 1, $false, "78×40=", 2) | Out-Null
$d.Content.Find.Execute("66×30=", $true, $false, $false, $false, $false, $true, 1, $false, "58×59=", 2) | Out-Null
$d.Content.Find.Execute("22×26=", $true, $false, $false, $false, $false, $true, 1, $false, "52×11=", 2) | Out-Null
$d.Content.Find.Execute("70×85=", $true, $false, $false, $false, $false, $true, 1, $false, "16×50=", 2) | Out-Null
$d.Content.Find.Execute("97×61=", $true, $false, $false, $false, $false, $true, 1, $false, "55×13=", 2) | Out-Null
$d.Content.Find.Execute("94×56=", $true, $false, $false, $false, $false, $true, 1, $false, "27×27=", 2) | Out-Null
$d.Content.Find.Execute("62×26=", $true, $false, $false, $false, $false, $true, 1, $false, "87×66=", 2) | Out-Null
$d.Content.Find.Execute("56×90=", $true, $false, $false, $false, $false, $true, 1, $false, "45×73=", 2) | Out-Null
$d.Content.Find.Execute("77×35=", $true, $false, $false, $false, $false, $true, 1, $false, "85×20=", 2) | Out-Null
$d.Content.Find.Execute("54×65=", $true, $false, $false, $false, $false, $true, 1, $false, "14×55=", 2) | Out-Null
$d.Content.Find.Execute("15×87=", $true, $false, $false, $false, $false, $true, 1, $false, "97×59=", 2) | Out-Null
$d.Content.Find.Execute("86×67=", $true, $false, $false, $false, $false, $true, 1, $false, "28×18=", 2) | Out-Null
$d.Content.Find.Execute("40×89=", $true, $false, $false, $false, $false, $true, 1, $false, "64×86=", 2) | Out-Null
$d.Content.Find.Execute("96×48=", $true, $false, $false, $false, $false, $true, 1, $false, "17×72=", 2) | Out-Null
$d.Content.Find.Execute("88×17=", $true, $false, $false, $false, $false, $true, 1, $false, "36×49=", 2) | Out-Null
